$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the label for the first UnitTest entry (row 24, column E)
$ws.Range("E24").Value = "WebServer - Plugin Manager, UnitTest - 1"

# Fix the Karin entry on row 23: hours reduced from 4 to 2.5
$ws.Range("D23").Value = 2.5

# Fill in the previously empty rows 25 and 26 with new entries
# Copy the date-column number format from row 24 so rows 25/26 match the
# existing "A" column look (center aligned, date number format) instead of
# the plain placeholder style the empty cells had before.
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A25:A26").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A25").Value = 41611
$ws.Range("B25").Value = "Teresa"
$ws.Range("D25").Value = 2.5
$ws.Range("E26").Value = "UnitTest - 2, 3, 4"
$ws.Range("E25").Value = "StaticFile - Fertig, SensorCloud - XML"

$ws.Range("A26").Value = 41611
$ws.Range("B26").Value = "Karin"
$ws.Range("D26").Value = 1.5

# Move the active cell selection to E29 (next empty row), matching the author's cursor move
$ws.Range("E29").Select()
